# daily auto push: 2026-02-15 05:02 UTC
# Insert a new log row for 2026/02/15 at row 798 (shifting the existing
# 798:839 block down to 799:840) on the only worksheet in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push everything from row 798 down by one row.
$ws.Rows.Item(798).Insert()

# Force column A to be read as literal text (not auto-parsed into a date
# serial) so it matches the rest of the "date" column, then populate the
# new row's four cells.
$ws.Range("A798").NumberFormat = "@"
$ws.Range("A798").Value = "2026/02/15"
$ws.Range("B798").Value = "日"
$ws.Range("C798").Value = 13
$ws.Range("D798").Value = 25

# Drop the number-format override so the new cell carries no style index,
# matching the plain (unstyled) data cells around it.
$ws.Range("A798").ClearFormats()
